$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull in the login URL for the HRMS demo site into a new cell next to the
# existing login data, as a live hyperlink.
$ws.Range("G1").Value = "https://hrmsdemo.onpassive.com"
$ws.Hyperlinks.Add($ws.Range("G1"), "https://hrmsdemo.onpassive.com")

# Match the existing hyperlink cell styling used by column B.
$ws.Range("G1").Style = "Hyperlink"

# Leave the new cell selected, matching the author's final selection.
$ws.Range("G1").Select()
